# mapping_Hungary_IND.xlsx — simplify steel description (remove RME)
#
# The author edited cell B2 (the long pipe-delimited material description)
# to drop the "/RME" token from the "24% S/LFM+CDH/RME/H:1" line, then
# (as a side effect of editing in Excel) the cell ended up with word-wrap
# turned on and an enlarged row height, and the sheet's last selection was
# left on B2:B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")

# 1. Remove "/RME" from the steel description text.
$text = $cell.Value2
$newText = $text.Replace("24% S/LFM+CDH/RME/H:1", "24% S/LFM+CDH/H:1")
$cell.Value = $newText

# 2. Turn on word-wrap for the cell (new cell style w/ wrapText alignment).
$cell.WrapText = $true

# 3. Expand row 2 to the (near-maximum) height seen after the edit.
$ws.Rows.Item(2).RowHeight = 409.6

# 4. Leave the selection spanning B2:B12 (as left by the author's last edit).
$ws.Range("B2:B12").Select() | Out-Null
